# Auto-generated edit script: updates cryptos list price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.790.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.053.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '525.08'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.42%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +4.91%  '
$ws.Range("E9").Value = '  +4.54%  '
$ws.Range("E10").Value = '  +8.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.371'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.26%  '
$ws.Range("E12").Value = '  +2.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.572.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.94'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000171'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +17.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.795.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("E17").Value = '  +7.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.053.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.57%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '340.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.502'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.73%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.97'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.00%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.173'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.87%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0₃0971'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.57%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.06%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.40'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +10.18%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.22%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.26%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.13'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.07%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.75'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.86%  '
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '156.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.98'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.81%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.33'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.65%  '
$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.85'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.83%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0708'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.26%  '
$ws.Range("B39").Value = 'RenzoRestakedETH'
$ws.Range("C39").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.087.86'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.42%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.20%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.10%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.45%  '
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.331.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.83%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.662'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.92%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.36%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.06%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0247'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.69%  '
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.06'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.90%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.61%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0896'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.06%  '
